$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 274.22
$ws.Range("I15").Value = 274.22
$ws.Range("K15").Value = 822.6600000000001
$ws.Range("M15").Value = -653.6600000000001
$ws.Range("H55").Value = 72.85714
$ws.Range("I55").Value = 68
$ws.Range("J55").Value = 85
$ws.Range("K55").Value = 68
$ws.Range("L55").Value = 85
$ws.Range("M55").Value = 146
$ws.Range("N55").Value = -513
$ws.Range("H109").Value = 79800
$ws.Range("J109").Value = 79800
$ws.Range("L109").Value = 79800
$ws.Range("N109").Value = -82574
$ws.Range("H132").Value = 2985
$ws.Range("I132").Value = 3218.5293
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 9655.5879
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -7125.5879
$ws.Range("N132").Value = -8060
$ws.Range("H135").Value = 37500300
$ws.Range("I135").Value = 13513796
$ws.Range("J135").Value = 333333860
$ws.Range("K135").Value = 121624164
$ws.Range("L135").Value = 3000004740
$ws.Range("M135").Value = -121621629
$ws.Range("N135").Value = -3000009810
$ws.Range("H137").Value = 1831.1666
$ws.Range("I137").Value = 1304.0731
$ws.Range("J137").Value = 2968.5789
$ws.Range("K137").Value = 3912.2193
$ws.Range("L137").Value = 8905.736699999999
$ws.Range("M137").Value = -1362.2193
$ws.Range("N137").Value = -14005.7367
$ws.Range("H138").Value = 3138.7595
$ws.Range("I138").Value = 1859.6
$ws.Range("J138").Value = 3730.963
$ws.Range("K138").Value = 5578.799999999999
$ws.Range("L138").Value = 11192.889
$ws.Range("M138").Value = -438.7999999999993
$ws.Range("N138").Value = -21472.889

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H53").Value = 800
$ws.Range("I53").Value = 800
$ws.Range("K53").Value = 800
$ws.Range("M53").Value = -118
$ws.Range("H61").Value = 7688.68
$ws.Range("I61").Value = 5930.1143
$ws.Range("K61").Value = 5930.1143
$ws.Range("M61").Value = -5718.1143
$ws.Range("H74").Value = 4839.3335
$ws.Range("I74").Value = 2265.4285
$ws.Range("J74").Value = 9343.666999999999
$ws.Range("K74").Value = 2265.4285
$ws.Range("L74").Value = 9343.666999999999
$ws.Range("M74").Value = -1391.4285
$ws.Range("N74").Value = -11091.667
$ws.Range("H77").Value = 4839.3335
$ws.Range("I77").Value = 2265.4285
$ws.Range("J77").Value = 9343.666999999999
$ws.Range("K77").Value = 11327.1425
$ws.Range("L77").Value = 46718.335
$ws.Range("M77").Value = -6959.1425
$ws.Range("N77").Value = -55454.335
$ws.Range("H136").Value = 7688.68
$ws.Range("I136").Value = 5930.1143
$ws.Range("K136").Value = 17790.3429
$ws.Range("M136").Value = -15240.3429

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 26558.354
$ws.Range("J112").Value = 26558.354
$ws.Range("L112").Value = 26558.354
$ws.Range("N112").Value = -29512.354

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4487.2856
$ws.Range("I31").Value = 4632.2
$ws.Range("J31").Value = 4125
$ws.Range("K31").Value = 4632.2
$ws.Range("L31").Value = 4125
$ws.Range("M31").Value = -4337.2
$ws.Range("N31").Value = -4715
$ws.Range("H34").Value = 4487.2856
$ws.Range("I34").Value = 4632.2
$ws.Range("J34").Value = 4125
$ws.Range("K34").Value = 4632.2
$ws.Range("L34").Value = 4125
$ws.Range("M34").Value = -4430.2
$ws.Range("N34").Value = -4529
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 10000338
$ws.Range("I18").Value = 11111376
$ws.Range("K18").Value = 33334128
$ws.Range("M18").Value = -33333959
$ws.Range("H25").Value = 2651
$ws.Range("J25").Value = 2651
$ws.Range("L25").Value = 7953
$ws.Range("N25").Value = -8291
$ws.Range("H30").Value = 2651
$ws.Range("J30").Value = 2651
$ws.Range("L30").Value = 7953
$ws.Range("N30").Value = -8157
$ws.Range("H58").Value = 2842.8572
$ws.Range("I58").Value = 266.66666
$ws.Range("J58").Value = 3041.0256
$ws.Range("K58").Value = 799.9999799999999
$ws.Range("L58").Value = 9123.076799999999
$ws.Range("M58").Value = -671.9999799999999
$ws.Range("N58").Value = -9379.076799999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8667.888999999999
$ws.Range("I80").Value = 26002.5
$ws.Range("J80").Value = 3715.1428
$ws.Range("K80").Value = 26002.5
$ws.Range("L80").Value = 3715.1428
$ws.Range("M80").Value = -25004.5
$ws.Range("N80").Value = -5711.1428
$ws.Range("H83").Value = 8667.888999999999
$ws.Range("I83").Value = 26002.5
$ws.Range("J83").Value = 3715.1428
$ws.Range("K83").Value = 130012.5
$ws.Range("L83").Value = 18575.714
$ws.Range("M83").Value = -125020.5
$ws.Range("N83").Value = -28559.714
$ws.Range("H112").Value = 49845
$ws.Range("J112").Value = 49845
$ws.Range("L112").Value = 49845
$ws.Range("N112").Value = -52061

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1967.1177
$ws.Range("I16").Value = 1185.4445
$ws.Range("J16").Value = 2846.5
$ws.Range("K16").Value = 1185.4445
$ws.Range("L16").Value = 2846.5
$ws.Range("M16").Value = -1015.4445
$ws.Range("N16").Value = -3186.5
$ws.Range("H61").Value = 483420.34
$ws.Range("I61").Value = 9497.076999999999
$ws.Range("J61").Value = 1253545.6
$ws.Range("K61").Value = 9497.076999999999
$ws.Range("L61").Value = 1253545.6
$ws.Range("M61").Value = -9295.076999999999
$ws.Range("N61").Value = -1253949.6
$ws.Range("H82").Value = 1900
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -3522
$ws.Range("H85").Value = 1900
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -5296
$ws.Range("H113").Value = 483420.34
$ws.Range("I113").Value = 9497.076999999999
$ws.Range("J113").Value = 1253545.6
$ws.Range("K113").Value = 9497.076999999999
$ws.Range("L113").Value = 1253545.6
$ws.Range("M113").Value = -7327.076999999999
$ws.Range("N113").Value = -1257885.6
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H136").Value = 3289.8513
$ws.Range("I136").Value = 2033.28
$ws.Range("J136").Value = 5907.7085
$ws.Range("K136").Value = 6099.84
$ws.Range("L136").Value = 17723.1255
$ws.Range("M136").Value = -3549.84
$ws.Range("N136").Value = -22823.1255

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H132").Value = 2917.9836
$ws.Range("I132").Value = 3087.205
$ws.Range("J132").Value = 2618
$ws.Range("K132").Value = 9261.615
$ws.Range("L132").Value = 7854
$ws.Range("M132").Value = -6731.615
$ws.Range("N132").Value = -12914
